# ------------------------------------------------------------------
# Edit script: syncs attendance_reports workbook per commit diff
#  1. Updates 19 rows on the "Summary" sheet with recalculated
#     attendance stats (G,I,N,O,X,Y columns) reflecting one extra
#     recorded BIOCHEMISTRY LAB/CBL session, and flips F column from
#     "Fail" to "High Risk" (w/ matching fill colour) for 4 students.
#  2. Appends 19 new attendance rows (576-594) to the "Attendance"
#     sheet for a BIOCHEMISTRY LAB/CBL session on 19/11/2025, widens
#     the Subject columns (F,H) and refreshes AutoFilter / the
#     hidden _FilterDatabase defined name to cover the new rows.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item("Summary")
$attendance = $wb.Worksheets.Item("Attendance")

# ---- helpers -------------------------------------------------------

function Set-TextValue($cell, $text) {
    # Force the cell to keep (or become) a text value even when the
    # text looks numeric / percentage-like, without permanently
    # altering the cell's number format.
    $origFormat = $cell.NumberFormat
    $cell.NumberFormat = "@"
    $cell.Value2 = $text
    $cell.NumberFormat = $origFormat
}

function Set-SummaryRow($row, $gPct, $iVal, $nVal, $oVal, $xVal, $yVal, $fStatus) {
    if ($fStatus -ne "") {
        # Status moved from "Fail" to "High Risk": swap the text and
        # recolour the cell to match the existing "High Risk" fill.
        $fCell = $summary.Range("F" + $row)
        $fCell.Interior.Color = 8158463
        $fCell.Value2 = $fStatus
    }

    Set-TextValue ($summary.Range("G" + $row)) $gPct
    $summary.Range("I" + $row).Value2 = $iVal
    $summary.Range("N" + $row).Value2 = $nVal
    $summary.Range("O" + $row).Value2 = $oVal
    $summary.Range("X" + $row).Value2 = $xVal
    $summary.Range("Y" + $row).Value2 = $yVal
}

function Add-AttendanceRow($row, $studentId, $name, $email, $date, $time) {
    $cA = $attendance.Cells.Item($row, 1)
    $cB = $attendance.Cells.Item($row, 2)
    $cC = $attendance.Cells.Item($row, 3)
    $cD = $attendance.Cells.Item($row, 4)
    $cE = $attendance.Cells.Item($row, 5)
    $cF = $attendance.Cells.Item($row, 6)
    $cG = $attendance.Cells.Item($row, 7)
    $cH = $attendance.Cells.Item($row, 8)
    $cI = $attendance.Cells.Item($row, 9)
    $cJ = $attendance.Cells.Item($row, 10)
    $cK = $attendance.Cells.Item($row, 11)

    Set-TextValue $cA $studentId
    Set-TextValue $cB $name
    Set-TextValue $cC "Year 2"
    Set-TextValue $cD "C1"
    Set-TextValue $cE $email
    Set-TextValue $cF "BIOCHEMISTRY LAB/CBL"
    Set-TextValue $cG "1"
    Set-TextValue $cH "BIOCHEMISTRY LAB/CBL"
    Set-TextValue $cI $date
    Set-TextValue $cJ $time
    Set-TextValue $cK "C1"

    # Cells added via COM pick up a non-default style the moment the
    # number format is touched; reset back to "Normal" so the row
    # matches the unstyled look of the rest of the attendance log.
    $cA.Style = "Normal"
    $cB.Style = "Normal"
    $cC.Style = "Normal"
    $cD.Style = "Normal"
    $cE.Style = "Normal"
    $cF.Style = "Normal"
    $cG.Style = "Normal"
    $cH.Style = "Normal"
    $cI.Style = "Normal"
    $cJ.Style = "Normal"
    $cK.Style = "Normal"
}

# ---- 1. Summary sheet recalculated rows -----------------------------

Set-SummaryRow 37 "6.9%" 21 2 9 1 1 ""
Set-SummaryRow 71 "17.2%" 18 5 6 1 1 "High Risk"
Set-SummaryRow 74 "17.2%" 18 5 6 1 1 "High Risk"
Set-SummaryRow 75 "17.2%" 18 5 6 1 1 "High Risk"
Set-SummaryRow 93 "13.8%" 19 4 7 1 1 ""
Set-SummaryRow 96 "13.8%" 19 4 7 1 1 ""
Set-SummaryRow 147 "20.7%" 17 6 5 1 1 ""
Set-SummaryRow 180 "3.4%" 22 1 10 1 1 ""
Set-SummaryRow 184 "20.7%" 17 6 5 1 1 ""
Set-SummaryRow 195 "17.2%" 18 5 6 1 1 "High Risk"
Set-SummaryRow 197 "13.8%" 19 4 7 1 1 ""
Set-SummaryRow 209 "13.8%" 19 4 7 1 1 ""
Set-SummaryRow 220 "10.3%" 20 3 8 1 1 ""
Set-SummaryRow 232 "20.7%" 17 6 5 1 1 ""
Set-SummaryRow 233 "27.6%" 15 8 3 1 1 ""
Set-SummaryRow 238 "20.7%" 17 6 5 1 1 ""
Set-SummaryRow 243 "20.7%" 17 6 5 1 1 ""
Set-SummaryRow 245 "20.7%" 17 6 5 1 1 ""
Set-SummaryRow 248 "10.3%" 20 3 8 1 1 ""

# ---- 2. Attendance sheet new BIOCHEMISTRY LAB/CBL rows --------------

Add-AttendanceRow 576 "221539" "تحريم شوكات مالك" "221539@med.asu.edu.eg" "19/11/2025" "11:28:04"
Add-AttendanceRow 577 "221755" "سعدية عاشق" "221755@med.asu.edu.eg" "19/11/2025" "11:28:16"
Add-AttendanceRow 578 "221833" "صفا محمود صايل صايل" "221833@med.asu.edu.eg" "19/11/2025" "11:28:49"
Add-AttendanceRow 579 "222058" "رفا السيد قسم الله السيد" "222058@med.asu.edu.eg" "19/11/2025" "11:31:08"
Add-AttendanceRow 580 "221904" "عائشه نور شيهو" "221904@med.asu.edu.eg" "19/11/2025" "11:31:54"
Add-AttendanceRow 581 "221948" "سانتينو اتيم شول دينق" "221948@med.asu.edu.eg" "19/11/2025" "11:32:33"
Add-AttendanceRow 582 "221714" "زينب عبد اللطيف بيبى فاروق" "221714@med.asu.edu.eg" "19/11/2025" "11:34:55"
Add-AttendanceRow 583 "221822" "سعادة يوسف عليو" "221822@med.asu.edu.eg" "19/11/2025" "11:35:26"
Add-AttendanceRow 584 "211620" "محمودول اسلام" "211620@med.asu.edu.eg" "19/11/2025" "11:35:49"
Add-AttendanceRow 585 "222076" "ابرار عبد الماجد عبد العزيز عثمان" "222076@med.asu.edu.eg" "19/11/2025" "11:36:26"
Add-AttendanceRow 586 "212442" "رميساء محى الدين الامين الطيب" "212442@med.asu.edu.eg" "19/11/2025" "11:36:42"
Add-AttendanceRow 587 "220304" "احمد الكامل محمد عبدون عثمان" "220304@med.asu.edu.eg" "19/11/2025" "11:37:00"
Add-AttendanceRow 588 "220967" "لارا حربي عبدالله الزيادات" "220967@med.asu.edu.eg" "19/11/2025" "11:37:18"
Add-AttendanceRow 589 "212543" "زينب سيف الدين محمد ادم" "212543@med.asu.edu.eg" "19/11/2025" "11:37:37"
Add-AttendanceRow 590 "222028" "هاجر عبد الحفيظ سيد صالح" "222028@med.asu.edu.eg" "19/11/2025" "11:37:49"
Add-AttendanceRow 591 "222003" "اسراء بدر الدين جعفر عثمان" "222003@med.asu.edu.eg" "19/11/2025" "11:37:54"
Add-AttendanceRow 592 "221000" "ابوبكر محمد قايد الثوابي" "221000@med.asu.edu.eg" "19/11/2025" "11:38:06"
Add-AttendanceRow 593 "222004" "احمد ايمن احمد بشير" "222004@med.asu.edu.eg" "19/11/2025" "11:38:18"
Add-AttendanceRow 594 "222053" "صباح سيف الدين عثمان اسحق" "222053@med.asu.edu.eg" "19/11/2025" "11:38:41"

# ---- 3. Widen Subject columns (F,H) on the Attendance sheet ---------

$attendance.Columns.Item(6).ColumnWidth = 21.166666666666668
$attendance.Columns.Item(8).ColumnWidth = 21.166666666666668

# ---- 4. Refresh AutoFilter range + hidden _FilterDatabase name ------

$attendance.AutoFilterMode = $False
$attendance.Range("A1:K594").AutoFilter()

$names = $wb.Names
for ($i = 1; $i -le $names.Count; $i++) {
    $n = $names.Item($i)
    if ($n.Name -eq "Attendance!_FilterDatabase") {
        $n.RefersTo = "='Attendance'!`$A`$1:`$K`$594"
    }
}
